$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182352423667908
$ws.Range("B1").Value = 2.189099311828613
$ws.Range("C1").Value = 4.439297199249268
$ws.Range("D1").Value = 2.705202102661133
$ws.Range("E1").Value = 1.219675660133362
